# Apply the "Agile Scrum Plan" sheet updates described by the commit:
#   "Replaced Fig [1] & Modified AI Excel Sheet"
#
# Changes:
#  - B8: "Created FSM" -> "Create OOP FSM"
#  - New retrospective/review columns (C = Yes/No, D = free-text comments)
#    filled in for rows 3-7, with legend labels added at J2/J3
#  - B9 picks up the "Good" (green) cell style
#  - Column widths adjusted for D (comments) and H (key), and a new
#    narrow hidden column I is introduced
#  - Selection moves to B9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agile Scrum Plan")

# --- Cell content -----------------------------------------------------

# Legend / key for the new Retrospect & Sprint review columns
$ws.Range("J2").Value = "Retrospect = Did we do it "
$ws.Range("J3").Value = "Sprint review = How well did we do it and what could be better"

# Retrospect column (C) - answered "Yes" for every logged sprint
$ws.Range("C3").Value = "Yes"
$ws.Range("C4").Value = "Yes"
$ws.Range("C5").Value = "Yes"
$ws.Range("C6").Value = "Yes"
$ws.Range("C7").Value = "Yes"

# Sprint review column (D) - free text summary for each sprint
$ws.Range("D3").Value = "We successfully researched FSM using different sources and Looked at multiple strategies of implementation"
$ws.Range("D4").Value = "We colour coded each week to correctly divide the different parts of the assignment plan, found no issues "
$ws.Range("D5").Value = "We created a lot of pseudocode for the FSM states, which we can now impliment into the project "
$ws.Range("D6").Value = "Created rules for each of the states and the varying types of tanks "
$ws.Range("D7").Value = "We started the introduction with the explanation of the FSM with also the project management "

# Project step renamed
$ws.Range("B8").Value = "Create OOP FSM"

# B9 gets the "Good" (green) built-in cell style
$ws.Range("B9").Style = "Good"

# --- Column sizing ------------------------------------------------------

# Column D (Sprint review comments) widened to fit the long text
$ws.Columns.Item(4).ColumnWidth = 97

# Column H (key) narrowed
$ws.Columns.Item(8).ColumnWidth = 9.83333333333333

# New column I, narrow and hidden
$ws.Columns.Item(9).ColumnWidth = 8.33333333333333
$ws.Columns.Item(9).Hidden = $true

# --- Selection ------------------------------------------------------

[void]$ws.Range("B9").Select()
